$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 6 ("grandes regiões e unidades da federação" header row),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(6).Delete()
